$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "G" column values (rows 5-18 and 20-26) to 51, leaving row 19 untouched.
$rows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,20,21,22,23,24,25,26)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = 51
}

# Update the active selection on the sheet to H22.
$ws.Range("H22").Select()
